$d = $word.ActiveDocument

# --- Paragraph containing the "m:userdoc 'zone1'" field -------------------
# Replace the field-code run sequence (fldChar begin / instrText* / fldChar end)
# with plain literal text runs: {  m  :userdoc 'zone1'  }
$p1 = $d.Paragraphs(2).Range
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:r><w:t>{</w:t></w:r>
  <w:r><w:t>m</w:t></w:r>
  <w:r><w:t>:userdoc 'zone1'</w:t></w:r>
  <w:r><w:t xml:space="preserve">}</w:t></w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@
$p1.InsertXML($xml1)

# --- Paragraph containing the "m:enduserdoc" field -------------------------
# Replace the field-code run sequence with plain literal text runs,
# keeping the bookmarkStart/bookmarkEnd ("_GoBack") between the two runs.
$p2 = $d.Paragraphs(4).Range
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:r><w:t>{m:</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:t xml:space="preserve">enduserdoc}</w:t></w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@
$p2.InsertXML($xml2)
